# Applies the Feb 2, 2023 05:47 UTC "Updated symbol list" data refresh:
# coin prices (column D), 1h volume deltas (column E), and several rows
# of coin name/link (columns B/C) that shifted position in the ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D/E hold numeric- and percent-looking text (e.g. "329.67", "6.04%")
# that must stay literal strings -- force text format before writing so the
# COM layer does not coerce them into real numbers and lose the formatting.
$updates = @(
    ,('D2', '329.67')
    ,('E2', '6.04%')
    ,('D3', '40.36')
    ,('E3', '7.12%')
    ,('D4', '5.268')
    ,('E4', '1.80%')
    ,('D5', '0.08119')
    ,('E5', '2.56%')
    ,('B6', 'KuCoinToken')
    ,('C6', 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs')
    ,('D6', '8.662')
    ,('E6', '4.51%')
    ,('B7', 'FTXToken')
    ,('C7', 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt')
    ,('D7', '1.909')
    ,('E7', '-0.73%')
    ,('B8', 'BTSEToken')
    ,('C8', 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse')
    ,('D8', '2.957')
    ,('E8', '-1.43%')
    ,('B9', 'MXToken')
    ,('C9', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx')
    ,('D9', '0.9354')
    ,('E9', '-0.35%')
    ,('B10', 'LiechtensteinCryptoassetsExchange')
    ,('C10', 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx')
    ,('D10', '0.1391')
    ,('E10', '26.36%')
    ,('B11', 'WazirX')
    ,('C11', 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx')
    ,('D11', '0.1967')
    ,('E11', '0.96%')
    ,('B12', 'MandalaExchangeToken')
    ,('C12', 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx')
    ,('D12', '0.09213')
    ,('E12', '1.52%')
    ,('B13', 'BitrueCoin')
    ,('C13', 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr')
    ,('D13', '0.03431')
    ,('E13', '3.80%')
    ,('B14', 'BitMartToken')
    ,('C14', 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx')
    ,('D14', '0.09559')
    ,('E14', '-0.59%')
    ,('B15', 'BitForexToken')
    ,('C15', 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf')
    ,('D15', '0.001391')
    ,('E15', '0.89%')
    ,('B16', 'TigerCash')
    ,('C16', 'https://coinranking.com/coin/6hIn06L2+tigercash-tch')
    ,('D16', '0.006012')
    ,('E16', '3.84%')
    ,('B17', 'LEO')
    ,('C17', 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo')
    ,('D17', '3.358')
    ,('E17', '-6.60%')
    ,('B18', 'GateToken')
    ,('C18', 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt')
    ,('D18', '4.529')
    ,('E18', '2.18%')
    ,('E19', '3.36%')
    ,('D20', '7.266')
    ,('E20', '12.95%')
    ,('D21', '0.1323')
    ,('E21', '3.21%')
    ,('E22', '1.82%')
    ,('E23', '0.87%')
    ,('E24', '-0.88%')
    ,('D25', '0.004358')
    ,('E25', '-5.57%')
    ,('E26', '-5.16%')
    ,('D27', '0.0003993')
    ,('E27', '-0.03%')
    ,('D39', '0.02523')
    ,('E39', '12.11%')
    ,('D40', '0.05224')
    ,('E40', '2.38%')
    ,('D41', '0.007720')
    ,('E41', '3.18%')
    ,('D42', '0.1432')
    ,('E42', '5.78%')
    ,('D43', '0.009045')
    ,('E43', '2.06%')
    ,('D44', '0.002171')
    ,('E44', '1.86%')
    ,('D45', '0.008966')
    ,('E45', '-3.76%')
    ,('D46', '0.00006636')
    ,('E46', '0.44%')
    ,('D47', '0.00000000751')
    ,('E47', '0.02%')
    ,('D48', '0.003342')
    ,('E49', '148.05%')
    ,('D50', '0.00002102')
    ,('E50', '0.02%')
    ,('D51', '0.0002001')
    ,('E51', '0.02%')
)

foreach ($pair in $updates) {
    $ref = $pair[0]
    $value = $pair[1]
    $range = $ws.Range($ref)
    $range.NumberFormat = "@"
    $range.Value = $value
}
